$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data as scraped on Sat Jul  6 04:53:02 UTC 2024.
# For D-column values that look numeric, force Text number format first so they
# are stored as literal strings (matching the original inline-string cell type)
# rather than being auto-converted to numbers by the Value setter.
$ws.Range("D2").Value = '56.144.48'
$ws.Range("E2").Value = '  +3.33%  '
$ws.Range("D3").Value = '2.966.50'
$ws.Range("E3").Value = '  +2.66%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '500.88'
$ws.Range("E5").Value = '  +5.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.21'
$ws.Range("E6").Value = '  +6.48%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.428'
$ws.Range("E8").Value = '  +6.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.33'
$ws.Range("E9").Value = '  +10.74%  '
$ws.Range("E10").Value = '  +8.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.350'
$ws.Range("E11").Value = '  +5.02%  '
$ws.Range("E12").Value = '  +3.14%  '
$ws.Range("D13").Value = '3.475.71'
$ws.Range("E13").Value = '  +2.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.20'
$ws.Range("E14").Value = '  +10.07%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000150'
$ws.Range("E15").Value = '  +11.51%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '56.178.56'
$ws.Range("E16").Value = '  +3.30%  '
$ws.Range("D17").Value = '2.969.21'
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.67'
$ws.Range("E18").Value = '  +9.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.27'
$ws.Range("E19").Value = '  +6.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.72'
$ws.Range("E20").Value = '  +8.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.56'
$ws.Range("E21").Value = '  +4.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("E23").Value = '  +4.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.74'
$ws.Range("E24").Value = '  +3.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.01'
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("E26").Value = '  +5.01%  '
$ws.Range("D27").Value = '0.0₃0881'
$ws.Range("E27").Value = '  +6.30%  '
$ws.Range("E28").Value = '  +3.64%  '
$ws.Range("E29").Value = '  +9.73%  '
$ws.Range("E30").Value = '  +2.49%  '
$ws.Range("E31").Value = '  +7.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.35'
$ws.Range("E32").Value = '  +6.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '158.47'
$ws.Range("E33").Value = '  +13.29%  '
$ws.Range("E34").Value = '  +4.47%  '
$ws.Range("E35").Value = '  +2.95%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.52'
$ws.Range("E36").Value = '  +1.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0670'
$ws.Range("E37").Value = '  +7.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.90'
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("D39").Value = '2.999.99'
$ws.Range("E39").Value = '  +2.92%  '
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.16'
$ws.Range("E41").Value = '  +3.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.637'
$ws.Range("E42").Value = '  +6.66%  '
$ws.Range("D43").Value = '2.235.07'
$ws.Range("E43").Value = '  +8.32%  '
$ws.Range("E44").Value = '  +5.85%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.970'
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("E46").Value = '  +2.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.91'
$ws.Range("E47").Value = '  +18.88%  '
$ws.Range("E48").Value = '  +7.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0233'
$ws.Range("E49").Value = '  +9.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.86'
$ws.Range("E50").Value = '  +3.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0864'
$ws.Range("E51").Value = '  +8.08%  '
